$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '26.173.66'
$ws.Range("D3").Value = '1.588.83'
$ws.Range("E3").Value = '  +0.02%  '
$ws.Range("E4").Value = '  -0.05%  '
Set-TextValue $ws.Range("D5") '211.69'
$ws.Range("E5").Value = '  +0.75%  '
$ws.Range("E6").Value = '  -1.09%  '
$ws.Range("E9").Value = '  -1.04%  '
Set-TextValue $ws.Range("D10") '18.96'
$ws.Range("E10").Value = '  -2.36%  '
$ws.Range("E11").Value = '  -0.12%  '
$ws.Range("D12").Value = '1.813.33'
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("D13").Value = '1.589.06'
$ws.Range("E13").Value = '  -0.63%  '
$ws.Range("E14").Value = '  -1.48%  '
$ws.Range("E15").Value = '  -1.92%  '
$ws.Range("D17").Value = '26.178.31'
$ws.Range("E17").Value = '  -0.51%  '
$ws.Range("D18").Value = '0.0₃0722'
$ws.Range("E18").Value = '  -0.89%  '
Set-TextValue $ws.Range("D19") '214.03'
$ws.Range("E19").Value = '  +1.51%  '
$ws.Range("E20").Value = '  -1.80%  '
$ws.Range("E21").Value = '  -0.08%  '
$ws.Range("E22").Value = '  -0.85%  '
$ws.Range("E23").Value = '  +0.59%  '
$ws.Range("E24").Value = '  -1.55%  '
Set-TextValue $ws.Range("D25") '144.66'
$ws.Range("E25").Value = '  +0.09%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("E27").Value = '  -1.35%  '
Set-TextValue $ws.Range("D28") '0.111'
$ws.Range("E28").Value = '  -1.35%  '
$ws.Range("E30").Value = '  -2.51%  '
$ws.Range("E31").Value = '  -0.01%  '
Set-TextValue $ws.Range("D32") '3.15'
$ws.Range("E32").Value = '  -1.89%  '
$ws.Range("D33").Value = '1.419.29'
$ws.Range("E33").Value = '  +7.99%  '
Set-TextValue $ws.Range("D34") '2.94'
$ws.Range("E34").Value = '  -1.84%  '
Set-TextValue $ws.Range("D35") '2.43'
$ws.Range("E35").Value = '  -0.78%  '
$ws.Range("E36").Value = '  -1.01%  '
Set-TextValue $ws.Range("D37") '0.584'
$ws.Range("E37").Value = '  -4.42%  '
$ws.Range("E38").Value = '  -1.82%  '
Set-TextValue $ws.Range("D39") '0.822'
$ws.Range("E39").Value = '  +1.98%  '
Set-TextValue $ws.Range("D40") '5.87'
$ws.Range("E40").Value = '  +4.65%  '
$ws.Range("E41").Value = '  -0.04%  '
Set-TextValue $ws.Range("D42") '0.965'
$ws.Range("E42").Value = '  -11.41%  '
$ws.Range("E43").Value = '  +0.37%  '
Set-TextValue $ws.Range("D44") '0.763'
$ws.Range("E44").Value = '  -0.59%  '
$ws.Range("D45").Value = '1.724.96'
$ws.Range("E45").Value = '  +0.00%  '
Set-TextValue $ws.Range("D46") '60.95'
$ws.Range("E46").Value = '  -2.08%  '
Set-TextValue $ws.Range("D47") '86.86'
$ws.Range("E47").Value = '  -0.82%  '
$ws.Range("E48").Value = '  -0.78%  '
$ws.Range("E49").Value = '  -0.88%  '
Set-TextValue $ws.Range("D50") '0.0959'
$ws.Range("E50").Value = '  -2.16%  '
Set-TextValue $ws.Range("D51") '0.999'
$ws.Range("E51").Value = '  -0.10%  '
